$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" header on "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename "Requested quantity" header on "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the current last sheet (Monthly Trend) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match sheet-level layout settings used by the other sheets in the workbook
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Reuse the same header/date cell formatting already used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A15").PasteSpecial(-4122)  # xlPasteFormats

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(44934.99999999999, 43, 11.83717235658957, 75.07254196586796),
    @(44969.99999999999, 43, 9.776996558810332, 74.23945611086972),
    @(44976.99999999999, 42, 10.27696198512871, 76.67918119414209),
    @(44983.99999999999, 42, 10.76868450455415, 76.66056760413589),
    @(44990.99999999999, 42, 10.24685076004642, 73.50415715327942),
    @(45004.99999999999, 42, 10.86089718774731, 74.13472552364051),
    @(45011.99999999999, 42, 7.866497626208765, 73.70394203019313),
    @(45018.99999999999, 42, 9.804358707137718, 74.4602790398468),
    @(45025.99999999999, 42, 11.63107725523667, 72.44740588558767),
    @(45032.99999999999, 42, 7.849383405773062, 74.45606157487582),
    @(45039.99999999999, 42, 10.6937839240331, 76.08349905960756),
    @(45046.99999999999, 42, 7.769114877289821, 73.04012678898614),
    @(45053.99999999999, 42, 9.916702353994433, 75.14896641573797),
    @(45060.99999999999, 41, 9.778710928972178, 72.68087989771188)
)

$row = 2
foreach ($rec in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
